$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.025074769481775
$ws.Range("D2").Value = 1.029645892457879
$ws.Range("E2").Value = 1.048426207647521
$ws.Range("F2").Value = 1.052710654235002
$ws.Range("I2").Value = 1.033204400483247
$ws.Range("J2").Value = 1.030245821730994
$ws.Range("K2").Value = 1.032459200633997
$ws.Range("L2").Value = 1.051186085976613
$ws.Range("M2").Value = 1.055458632899228
$ws.Range("N2").Value = 1.014145006689361
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.02582287612453
$ws.Range("D3").Value = 1.030185288529728
$ws.Range("E3").Value = 1.049607812508344
$ws.Range("F3").Value = 1.053942035605907
$ws.Range("I3").Value = 1.03335362063743
$ws.Range("J3").Value = 1.030634147218604
$ws.Range("K3").Value = 1.03280750923893
$ws.Range("L3").Value = 1.0521787776471
$ws.Range("M3").Value = 1.056501841028601
$ws.Range("N3").Value = 1.014272898768009
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.026307101972126
$ws.Range("D4").Value = 1.030534332644254
$ws.Range("E4").Value = 1.050373368525781
$ws.Range("F4").Value = 1.054739769525442
$ws.Range("I4").Value = 1.033448805515501
$ws.Range("J4").Value = 1.03088488746025
$ws.Range("K4").Value = 1.033032198492007
$ws.Range("L4").Value = 1.052821519744478
$ws.Range("M4").Value = 1.057177249074028
$ws.Range("N4").Value = 1.014355474647787
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.02651070513779
$ws.Range("D5").Value = 1.030681073698293
$ws.Range("E5").Value = 1.050695442955047
$ws.Range("F5").Value = 1.055075364121302
$ws.Range("I5").Value = 1.03348849243943
$ws.Range("J5").Value = 1.030990170258391
$ws.Range("K5").Value = 1.03312649169322
$ws.Range("L5").Value = 1.053091825579523
$ws.Range("M5").Value = 1.057461282410824
$ws.Range("N5").Value = 1.01439014639881
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.026544892991989
$ws.Range("D6").Value = 1.030705712291169
$ws.Range("E6").Value = 1.050749534442044
$ws.Range("F6").Value = 1.055131725267197
$ws.Range("I6").Value = 1.033495136747623
$ws.Range("J6").Value = 1.031007840146047
$ws.Range("K6").Value = 1.033142314162453
$ws.Range("L6").Value = 1.053137216818687
$ws.Range("M6").Value = 1.057508978246505
$ws.Range("N6").Value = 1.014395965398417
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.026309822393429
$ws.Range("D7").Value = 1.030536293395981
$ws.Range("E7").Value = 1.050377671177439
$ws.Range("F7").Value = 1.054744252863451
$ws.Range("I7").Value = 1.033449337106939
$ws.Range("J7").Value = 1.030886294759884
$ws.Range("K7").Value = 1.033033459096591
$ws.Range("L7").Value = 1.052825131202533
$ws.Range("M7").Value = 1.057181043980434
$ws.Range("N7").Value = 1.014355938103386
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.025327563500636
$ws.Range("D8").Value = 1.029828179041842
$ws.Range("E8").Value = 1.0488253341346
$ws.Range("F8").Value = 1.053126609284014
$ws.Range("I8").Value = 1.033255113498542
$ws.Range("J8").Value = 1.030377167789422
$ws.Range("K8").Value = 1.032577055339087
$ws.Range("L8").Value = 1.051521487252904
$ws.Range("M8").Value = 1.055811110859174
$ws.Range("N8").Value = 1.014188265247971
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.023597918595111
$ws.Range("D9").Value = 1.02858060496783
$ws.Range("E9").Value = 1.046097395926719
$ws.Range("F9").Value = 1.050283358395037
$ws.Range("I9").Value = 1.032902400335583
$ws.Range("J9").Value = 1.029475986473076
$ws.Range("K9").Value = 1.031767577252818
$ws.Range("L9").Value = 1.049227382747422
$ws.Range("M9").Value = 1.053400035362274
$ws.Range("N9").Value = 1.013891449134885
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.022445730546922
$ws.Range("D10").Value = 1.027749126812908
$ws.Range("E10").Value = 1.04428377154272
$ws.Range("F10").Value = 1.048392703849908
$ws.Range("I10").Value = 1.03266026185233
$ws.Range("J10").Value = 1.028872542938583
$ws.Range("K10").Value = 1.031224466141091
$ws.Range("L10").Value = 1.047700031938739
$ws.Range("M10").Value = 1.051794601335236
$ws.Range("N10").Value = 1.013692680193411
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.021947054193287
$ws.Range("D11").Value = 1.027389163779307
$ws.Range("E11").Value = 1.043499631646111
$ws.Range("F11").Value = 1.047575172746192
$ws.Range("I11").Value = 1.032553762219434
$ws.Range("J11").Value = 1.02861062703743
$ws.Range("K11").Value = 1.030988485447699
$ws.Range("L11").Value = 1.047039153678492
$ws.Range("M11").Value = 1.051099889034188
$ws.Range("N11").Value = 1.013606403663556
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.021761859396075
$ws.Range("D12").Value = 1.027255469959158
$ws.Range("E12").Value = 1.043208542501455
$ws.Range("F12").Value = 1.047271675149257
$ws.Range("I12").Value = 1.032513955961265
$ws.Range("J12").Value = 1.0285132474559
$ws.Range("K12").Value = 1.030900711032375
$ws.Range("L12").Value = 1.04679374484904
$ws.Range("M12").Value = 1.050841909298435
$ws.Range("N12").Value = 1.01357432574988
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.021801582665293
$ws.Range("D13").Value = 1.027284147147362
$ws.Range("E13").Value = 1.043270974187848
$ws.Range("F13").Value = 1.047336768737832
$ws.Range("I13").Value = 1.03252250573058
$ws.Range("J13").Value = 1.0285341398711
$ws.Range("K13").Value = 1.030919544387137
$ws.Range("L13").Value = 1.046846382681573
$ws.Range("M13").Value = 1.050897243787096
$ws.Range("N13").Value = 1.01358120796753
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.021931745194816
$ws.Range("D14").Value = 1.027378112342732
$ws.Range("E14").Value = 1.043475566553828
$ws.Range("F14").Value = 1.04755008207879
$ws.Range("I14").Value = 1.032550476867945
$ws.Range("J14").Value = 1.028602579489501
$ws.Range("K14").Value = 1.030981232443104
$ws.Range("L14").Value = 1.047018866667256
$ws.Range("M14").Value = 1.051078562977498
$ws.Range("N14").Value = 1.013603752723112
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.02201194741284
$ws.Range("D15").Value = 1.027436009101619
$ws.Range("E15").Value = 1.043601645875573
$ws.Range("F15").Value = 1.047681533957059
$ws.Range("I15").Value = 1.032567678025217
$ws.Range("J15").Value = 1.028644735186993
$ws.Range("K15").Value = 1.031019224529523
$ws.Range("L15").Value = 1.047125149106818
$ws.Range("M15").Value = 1.051190288625597
$ws.Range("N15").Value = 1.013617639196975
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.022478830995745
$ws.Range("D16").Value = 1.027773018036881
$ws.Range("E16").Value = 1.044335836948715
$ws.Range("F16").Value = 1.048446984551252
$ws.Range("I16").Value = 1.032667295130682
$ws.Range("J16").Value = 1.028889912432323
$ws.Range("K16").Value = 1.031240110423455
$ws.Range("L16").Value = 1.047743902207934
$ws.Range("M16").Value = 1.051840716571222
$ws.Range("N16").Value = 1.01369840172201
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.022771756978169
$ws.Range("D17").Value = 1.027984435369163
$ws.Range("E17").Value = 1.044796688763876
$ws.Range("F17").Value = 1.048927434852342
$ws.Range("I17").Value = 1.032729340463503
$ws.Range("J17").Value = 1.029043540070084
$ws.Range("K17").Value = 1.031378450166612
$ws.Range("L17").Value = 1.048132156394381
$ws.Range("M17").Value = 1.052248833326534
$ws.Range("N17").Value = 1.013749006388085
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.022942637714576
$ws.Range("D18").Value = 1.028107758410826
$ws.Range("E18").Value = 1.045065608901034
$ws.Range("F18").Value = 1.049207782917381
$ws.Range("I18").Value = 1.032765371011743
$ws.Range("J18").Value = 1.029133088505348
$ws.Range("K18").Value = 1.031459063177375
$ws.Range("L18").Value = 1.048358664331889
$ws.Range("M18").Value = 1.052486924491589
$ws.Range("N18").Value = 1.013778503117496
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.023000907311585
$ws.Range("D19").Value = 1.028149809501149
$ws.Range("E19").Value = 1.045157322917383
$ws.Range("F19").Value = 1.049303392960798
$ws.Range("I19").Value = 1.032777629433131
$ws.Range("J19").Value = 1.029163612008083
$ws.Range("K19").Value = 1.031486536832289
$ws.Range("L19").Value = 1.048435905476727
$ws.Range("M19").Value = 1.052568114829363
$ws.Range("N19").Value = 1.013788557318143
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.022740326509872
$ws.Range("D20").Value = 1.027961751566247
$ws.Range("E20").Value = 1.04474723203485
$ws.Range("F20").Value = 1.048875875756895
$ws.Range("I20").Value = 1.03272270007387
$ws.Range("J20").Value = 1.029027063474409
$ws.Range("K20").Value = 1.031363615704819
$ws.Range("L20").Value = 1.048090495656975
$ws.Range("M20").Value = 1.052205041757885
$ws.Range("N20").Value = 1.013743579064718
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.021893414574098
$ws.Range("D21").Value = 1.02735044157756
$ws.Range("E21").Value = 1.043415314366824
$ws.Range("F21").Value = 1.047487261932784
$ws.Range("I21").Value = 1.032542246893582
$ws.Range("J21").Value = 1.028582428269706
$ws.Range("K21").Value = 1.030963070171352
$ws.Range("L21").Value = 1.046968072482889
$ws.Range("M21").Value = 1.051025167139693
$ws.Range("N21").Value = 1.013597114706833
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.021361135293202
$ws.Range("D22").Value = 1.026966159698461
$ws.Range("E22").Value = 1.042578898926667
$ws.Range("F22").Value = 1.046615167730862
$ws.Range("I22").Value = 1.032427356876193
$ws.Range("J22").Value = 1.028302334592856
$ws.Range("K22").Value = 1.030710533586463
$ws.Range("L22").Value = 1.046262769745224
$ws.Range("M22").Value = 1.050283722558564
$ws.Range("N22").Value = 1.013504847732656
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.021643286424583
$ws.Range("D23").Value = 1.02716986734841
$ws.Range("E23").Value = 1.043022202676787
$ws.Range("F23").Value = 1.047077388455423
$ws.Range("I23").Value = 1.032488397762018
$ws.Range("J23").Value = 1.028450867847758
$ws.Range("K23").Value = 1.030844473805337
$ws.Range("L23").Value = 1.046636625456449
$ws.Range("M23").Value = 1.050676739594675
$ws.Range("N23").Value = 1.013553777063563
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.022754528517867
$ws.Range("D24").Value = 1.027972001379277
$ws.Range("E24").Value = 1.04476957905223
$ws.Range("F24").Value = 1.048899172754091
$ws.Range("I24").Value = 1.032725701072836
$ws.Range("J24").Value = 1.029034508724148
$ws.Range("K24").Value = 1.03137031900136
$ws.Range("L24").Value = 1.048109320209083
$ws.Range("M24").Value = 1.05222482914794
$ws.Range("N24").Value = 1.013746031500719
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.024044919542403
$ws.Range("D25").Value = 1.028903097704092
$ws.Range("E25").Value = 1.046801749934501
$ws.Range("F25").Value = 1.051017551000497
$ws.Range("I25").Value = 1.032994821270744
$ws.Range("J25").Value = 1.029709435716261
$ws.Range("K25").Value = 1.031977460918907
$ws.Range("L25").Value = 1.049820100563308
$ws.Range("M25").Value = 1.054023011475923
$ws.Range("N25").Value = 1.013968341782731
